# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.337.54"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.222.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.555"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.91"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.562.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.198.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.801"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.218.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.90%  "
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.05%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "146.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0754"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -4.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0288"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.748.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "78.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.179"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "66.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  -6.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.442.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.87%  "
